$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.547.03"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'3.054.92"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'386.56"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'102.99"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'36.83"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'3.548.53"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "'18.54"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "'3.033.70"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "'0.975"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "'10.69"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").Value = "'51.614.05"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'3.17"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "'12.43"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "'0.0₃0966"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'70.21"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'268.35"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'3.16"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").Value = "'8.23"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "'26.81"
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.170"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.25"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'34.73"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'50.02"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'0.0446"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").Value = "'0.294"
$ws.Range("E39").Value = "  +7.82%  "
$ws.Range("D40").Value = "'16.95"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "'125.32"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'3.74"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").Value = "'21.87"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "'2.034.10"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "'3.361.16"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("E51").Value = "  +6.09%  "
